# Exploratory data analysis edit on the Revenue sheet:
#  - insert two new columns (Average F&B, ARR) between "Marketing" and "Resturant Revenue"
#  - populate the new columns with per-row occupancy/ARR figures
#  - add forecast rows 22-25 (percentage achieved + projected revenue figures)
#  - fix a stray leading space in the "Local Rainy Season" label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revenue")

# Fix stray leading space on the "Local Rainy Season" label (col F)
$ws.Range("F1").Value = "Local Rainy Season"

# --- Insert the two new columns right before the old "Resturant Revenue" column (I) ---
$ws.Range("I1:J1").EntireColumn.Insert()

# New column headers
$ws.Range("I1").Value = "Average F&B"
$ws.Range("J1").Value = "ARR"

# Column widths for the two newly inserted columns
$ws.Columns.Item(9).ColumnWidth = 18.33
$ws.Columns.Item(10).ColumnWidth = 12.67

# --- Fill Average F&B (I) / ARR (J) values for every data row ---
$avgFB = @{
    2=15; 3=15; 4=15; 5=15; 6=15; 7=15;
    8=17; 9=17; 10=17; 11=17; 12=17; 13=17;
    14=19; 15=19; 16=19; 17=19; 18=19; 19=19; 20=19;
    21=20; 22=20; 23=20; 24=20; 25=20
}
$arr = @{
    2=75; 3=75; 4=75; 5=75; 6=75; 7=75;
    8=75; 9=77; 10=77; 11=77; 12=77; 13=77;
    14=79; 15=79; 16=79; 17=79; 18=79; 19=79; 20=79;
    21=80; 22=80; 23=80; 24=80; 25=80
}

foreach ($r in 2..25) {
    $ws.Cells.Item($r, 9).Value = $avgFB[$r]
    $ws.Cells.Item($r, 10).Value = $arr[$r]
}

# --- Holiday season category bump for the last two rows ---
$ws.Range("G24").Value = 5
$ws.Range("G25").Value = 5

# --- New forecast rows 22-25: % achieved (col B) ---
$ws.Range("B22").Value = 0.7
$ws.Range("B22").Style = "Percent"
$ws.Range("B23").Value = 0.4
$ws.Range("B23").Style = "Percent"
$ws.Range("B24").Value = 0.6
$ws.Range("B24").Style = "Percent"
$ws.Range("B25").Value = 0.8
$ws.Range("B25").Style = "Percent"

# --- Projected revenue figures (cols K:N), each row grown 7% off an earlier month ---
foreach ($pair in @(@(22,10), @(23,11), @(24,12), @(25,13))) {
    $destRow = $pair[0]
    $srcRow = $pair[1]
    foreach ($col in @("K", "L", "M", "N")) {
        $ws.Range("$col$destRow").Formula = "=$col$srcRow*1.07"
    }
}
